$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the condition string in A4: turn it into a comment line by
# prefixing it with '#' (matches the other comment rows in the sheet).
$cell = $ws.Range("A4")
$cell.Value2 = "#" + $cell.Value2

# Move/restore the active selection to A4 (was B7).
[void]$ws.Range("A4").Select()
